$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the FirstName/LastName header labels in B1/C1
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Shift the CNE identifiers in column A (rows 2-11) down by 1000010
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1000010
}

# Re-stamp the CNE / LastName / FirstName columns (A:C) with a fresh,
# uniform style across the header and all student rows (as re-emitted by
# the conversion tool when regenerating this sheet).
$ws.Range("A1:C1").Style = "Normal"
$ws.Range("A1:C11").WrapText = $false

# Update the active selection
$ws.Range("F7").Select()
